$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Octubre de 2020 a las 00:18"

# Colombia overtook Espana in total cases -> rows 8/9 swap which country is shown,
# Colombia gets the newer (higher) figures, Espana keeps its previous figures.
$ws.Range("A8").Value = "Colombia"
$ws.Range("B8").Value = 886179
$ws.Range("C8").Value = 8496
$ws.Range("D8").Value = 777658
$ws.Range("E8").Value = 81190
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 151
$ws.Range("H8").Value = 27331

$ws.Range("A9").Value = "España"
$ws.Range("B9").Value = 884381
$ws.Range("C9").Value = 5585
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 126
$ws.Range("H9").Value = 32688

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7825039
$ws.Range("C4").Value = 48176
$ws.Range("D4").Value = 5014922
$ws.Range("E4").Value = 2592534
$ws.Range("G4").Value = 802
$ws.Range("H4").Value = 217583

# Row 6 - Brasil
$ws.Range("B6").Value = 5028444
$ws.Range("C6").Value = 26087
$ws.Range("D6").Value = 4414564
$ws.Range("E6").Value = 464923
$ws.Range("G6").Value = 653
$ws.Range("H6").Value = 148957

# Row 13 - Sudafrica
$ws.Range("B13").Value = 686891
$ws.Range("C13").Value = 1736
$ws.Range("D13").Value = 618771
$ws.Range("E13").Value = 50712
$ws.Range("G13").Value = 160
$ws.Range("H13").Value = 17408

# Row 29 - Bolivia
$ws.Range("B29").Value = 175449
$ws.Range("C29").Value = 2326
$ws.Range("D29").Value = 147406
$ws.Range("E29").Value = 18487
$ws.Range("G29").Value = 15
$ws.Range("H29").Value = 9556

# Egipto overtook Oman -> rows 42/43 swap which country is shown,
# Egipto gets the newer (higher) figures, Oman keeps its previous figures.
$ws.Range("A42").Value = "Egipto"
$ws.Range("B42").Value = 104156
$ws.Range("C42").Value = 121
$ws.Range("D42").Value = 97524
$ws.Range("E42").Value = 615
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 7
$ws.Range("H42").Value = 6017

$ws.Range("A43").Value = "Oman"
$ws.Range("B43").Value = 104129
$ws.Range("C43").Value = 664
$ws.Range("D43").Value = 91731
$ws.Range("E43").Value = 11389
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 9
$ws.Range("H43").Value = 1009

# Row 57 - Bielorrusia
$ws.Range("B57").Value = 74422
$ws.Range("C57").Value = 490
$ws.Range("E57").Value = 4747

# Row 111 - Guayana Francesa
$ws.Range("B111").Value = 9219
$ws.Range("C111").Value = 100
$ws.Range("D111").Value = 7973
$ws.Range("E111").Value = 1117
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 129

# Row 145 - Somalia
$ws.Range("B145").Value = 3329
$ws.Range("C145").Value = 37
$ws.Range("D145").Value = 2165
$ws.Range("E145").Value = 1066
$ws.Range("G145").Value = 3
$ws.Range("H145").Value = 98

# Row 157 - Yemen
$ws.Range("B157").Value = 2241
$ws.Range("C157").Value = 19
$ws.Range("D157").Value = 1506
$ws.Range("E157").Value = 675
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 60

# Row 159 - Belice
$ws.Range("B159").Value = 2050
$ws.Range("C159").Value = 1
$ws.Range("D159").Value = 1329
